$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (date moved from 12-19 to 12-20)
$ws.Name = "Through 2021-12-20"

# Update header label for December row
$ws.Range("A13").Value = "December (through 12-20)"

# Update December row (row 13) values for each year column
$ws.Range("B13").Value = 27
$ws.Range("C13").Value = 64
$ws.Range("D13").Value = 77
$ws.Range("E13").Value = 44
$ws.Range("F13").Value = 33
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 138

# Update Total row (row 14) values for each year column
$ws.Range("B14").Value = 318
$ws.Range("C14").Value = 627
$ws.Range("D14").Value = 898
$ws.Range("E14").Value = 726
$ws.Range("F14").Value = 567
$ws.Range("G14").Value = 1357
$ws.Range("H14").Value = 1781
